$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.495.78"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "3.057.73"
$ws.Range("E3").Value = "  -2.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.30"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.42"
$ws.Range("E6").Value = "  +6.15%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  +3.20%  "
$ws.Range("D9").Value = "3.078.04"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.84"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.57"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "3.564.50"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("D18").Value = "63.432.13"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "3.069.40"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.45"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +5.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.69"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "445.77"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.31"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.111"
$ws.Range("E46").Value = "  +3.18%  "
$ws.Range("D47").Value = "2.807.47"
$ws.Range("E47").Value = "  -3.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.04"
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +4.04%  "
$ws.Range("E51").Value = "  +0.94%  "
